$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slurry pH")
$ws.Activate()

$ws.Range("A9").Value = 6.9
$ws.Range("A10").Select()
